# Re-process the metadata sheet with the newly curated dimensions.
# Column B (subespecie-ganaderia) moves from being a dimension to a measure,
# and column D (municipio-nombre) moves from being a measure to a dimension
# (classified like the other "refArea" dimension columns E/J), and the
# now-unused mapping file reference for subespecie-ganaderia is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B: subespecie-ganaderia -- dimension -> measure
$ws.Range("B2").Value = "iaest-measure:subespecie-ganaderia"
$ws.Range("B3").Value = "medida"
$ws.Range("B4").Value = "xsd:int"

# Column D: municipio-nombre -- measure -> dimension (refArea-like)
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# The mapping file for subespecie-ganaderia is no longer needed now that it
# is a measure instead of a curated dimension.
$ws.Range("B5").Value = ""
